$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 text with new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = $ws1.Range("A1").Value2
$text = $text -replace [regex]::Escape("1000 Bs = 15.15 = 61920.48 pesos"), "1000 Bs = 15.17 = 61883.31 pesos"
$text = $text -replace [regex]::Escape("61920.48 pesos = 15.07 = 981.09 Bs"), "61883.31 pesos = 15.06 = 959.87 Bs"
$ws1.Range("A1").Value2 = $text

# --- Update "tasas" sheet rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value2 = 65.90000000000001
$ws2.Range("O10").Value2 = 4078.11
$ws2.Range("O12").Value2 = 63.75
